{"js": "const replacements = [\n  [\"2025-02-26 Wednesday\", \"2025-02-27 Thursday\"],\n  [\"26\u00d784=2184\", \"83\u00d736=2988\"],\n  [\"21\u00d726=546\", \"39\u00d795=3705\"],\n  [\"91\u00d717=1547\", \"60\u00d751=3060\"],\n  [\"19\u00d774=1406\", \"83\u00d723=1909\"],\n  [\"42\u00d712=504\", \"94\u00d731=2914\"],\n  [\"47\u00d795=4465\", \"29\u00d740=1160\"],\n  [\"11\u00d792=1012\", \"45\u00d744=1980\"],\n  [\"80\u00d762=4960\", \"34\u00d795=3230\"],\n  [\"92\u00d732=2944\", \"22\u00d770=1540\"],\n  [\"30\u00d720=600\", \"14\u00d754=756\"],\n  [\"82\u00d751=4182\", \"77\u00d727=2079\"],\n  [\"20\u00d734=680\", \"25\u00d772=1800\"],\n  [\"79\u00d780=6320\", \"26\u00d750=1300\"],\n  [\"12\u00d741=492\", \"88\u00d745=3960\"],\n  [\"91\u00d761=5551\", \"54\u00d754=2916\"],\n  [\"65\u00d785=5525\", \"20\u00d750=1000\"],\n  [\"43\u00d729=1247\", \"13\u00d787=1131\"],\n  [\"30\u00d786=2580\", \"34\u00d795=3230\"],\n  [\"61\u00d758=3538\", \"76\u00d722=1672\"],\n  [\"75\u00d758=4350\", \"63\u00d742=2646\"],\n  [\"83\u00d773=6059\", \"46\u00d749=2254\"],\n  [\"94\u00d720=1880\", \"78\u00d719=1482\"],\n  [\"39\u00d769=2691\", \"58\u00d732=1856\"],\n  [\"70\u00d752=3640\", \"63\u00d795=5985\"],\n  [\"14\u00d789=1246\", \"56\u00d741=2296\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-26 Wednesday\", \"2025-02-27 Thursday\"),\n    @(\"26\u00d784=2184\", \"83\u00d736=2988\"),\n    @(\"21\u00d726=546\", \"39\u00d795=3705\"),\n    @(\"91\u00d717=1547\", \"60\u00d751=3060\"),\n    @(\"19\u00d774=1406\", \"83\u00d723=1909\"),\n    @(\"42\u00d712=504\", \"94\u00d731=2914\"),\n    @(\"47\u00d795=4465\", \"29\u00d740=1160\"),\n    @(\"11\u00d792=1012\", \"45\u00d744=1980\"),\n    @(\"80\u00d762=4960\", \"34\u00d795=3230\"),\n    @(\"92\u00d732=2944\", \"22\u00d770=1540\"),\n    @(\"30\u00d720=600\", \"14\u00d754=756\"),\n    @(\"82\u00d751=4182\", \"77\u00d727=2079\"),\n    @(\"20\u00d734=680\", \"25\u00d772=1800\"),\n    @(\"79\u00d780=6320\", \"26\u00d750=1300\"),\n    @(\"12\u00d741=492\", \"88\u00d745=3960\"),\n    @(\"91\u00d761=5551\", \"54\u00d754=2916\"),\n    @(\"65\u00d785=5525\", \"20\u00d750=1000\"),\n    @(\"43\u00d729=1247\", \"13\u00d787=1131\"),\n    @(\"30\u00d786=2580\", \"34\u00d795=3230\"),\n    @(\"61\u00d758=3538\", \"76\u00d722=1672\"),\n    @(\"75\u00d758=4350\", \"63\u00d742=2646\"),\n    @(\"83\u00d773=6059\", \"46\u00d749=2254\"),\n    @(\"94\u00d720=1880\", \"78\u00d719=1482\"),\n    @(\"39\u00d769=2691\", \"58\u00d732=1856\"),\n    @(\"70\u00d752=3640\", \"63\u00d795=5985\"),\n    @(\"14\u00d789=1246\", \"56\u00d741=2296\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute(\n        $oldText,  # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n    if (-not $result) {\n        throw \"Replacement failed for: $oldText\"\n    }\n}"}
